$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.931.96"
$ws.Range("E2").Value = "  +1.88%  "
$ws.Range("D3").Value = "2.246.15"
$ws.Range("E3").Value = "  +1.15%  "
$ws.Range("E4").Value = "  +0.02%  "
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "319.26"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  +0.44%  "
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.41"
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = "  +2.29%  "
$ws.Range("E7").Value = "  -1.05%  "
$ws.Range("E8").Value = "  +0.06%  "
$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.546"
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = "  -2.91%  "
$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.82"
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = "  -0.18%  "
$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0828"
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = "  +0.29%  "
$ws.Range("E12").Value = "  -1.58%  "
$ws.Range("E13").Value = "  -1.61%  "
$ws.Range("D14").Value = "2.587.84"
$ws.Range("E14").Value = "  +0.97%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$origStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.852"
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = "  -0.84%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$origStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.34"
$ws.Range("D16").Style = $origStyle
$ws.Range("E16").Value = "  +0.38%  "
$ws.Range("D17").Value = "2.249.74"
$ws.Range("E17").Value = "  +1.01%  "
$ws.Range("D18").Value = "43.798.98"
$ws.Range("E18").Value = "  +1.78%  "
$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "13.52"
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = "  -1.05%  "
$ws.Range("D20").Value = "0.0₃0977"
$ws.Range("E20").Value = "  +1.79%  "
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.44"
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = "  -1.60%  "
$origStyle = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "65.40"
$ws.Range("D22").Style = $origStyle
$ws.Range("E22").Value = "  +0.58%  "
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.09"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  -4.23%  "
$origStyle = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "233.76"
$ws.Range("D24").Style = $origStyle
$ws.Range("E24").Value = "  -0.81%  "
$ws.Range("E25").Value = "  -5.76%  "
$ws.Range("E26").Value = "  +0.20%  "
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.62"
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = "  +5.97%  "
$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "38.81"
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = "  +5.45%  "
$ws.Range("E29").Value = "  -1.38%  "
$ws.Range("E30").Value = "  -4.18%  "
$origStyle = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "159.63"
$ws.Range("D31").Style = $origStyle
$ws.Range("E31").Value = "  +1.32%  "
$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "20.10"
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = "  -0.26%  "
$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0841"
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = "  -2.40%  "
$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.67"
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = "  +0.11%  "
$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.113"
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = "  +9.48%  "
$origStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.06"
$ws.Range("D36").Style = $origStyle
$ws.Range("E36").Value = "  -8.01%  "
$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.94"
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = "  +5.38%  "
$ws.Range("E38").Value = "  -1.37%  "
$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.39"
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = "  +16.89%  "
$ws.Range("E40").Value = "  +0.35%  "
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.16"
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = "  -5.14%  "
$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0313"
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = "  -1.34%  "
$ws.Range("E43").Value = "  +0.14%  "
$ws.Range("D44").Value = "1.764.03"
$ws.Range("E44").Value = "  +1.12%  "
$origStyle = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "74.58"
$ws.Range("D45").Style = $origStyle
$ws.Range("E45").Value = "  +1.71%  "
$ws.Range("E46").Value = "  -3.10%  "
$ws.Range("B47").Value = "BitcoinSV"
$ws.Range("C47").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "81.02"
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = "  -2.52%  "
$ws.Range("B48").Value = "THORChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.15"
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = "  -1.88%  "
$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "103.65"
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = "  +1.28%  "
$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.67"
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = "  +4.55%  "
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "57.21"
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = "  -0.13%  "
